# Auto-generated edit script: apply numeric corrections to H:N columns
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# per the commit's scheduled-runner price-refresh diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 395.2
$ws.Range("I2").Value = 325.33334
$ws.Range("K2").Value = 325.33334
$ws.Range("M2").Value = -212.33334
# row 28
$ws.Range("H28").Value = 1446.5834
$ws.Range("I28").Value = 1494.125
$ws.Range("K28").Value = 1494.125
$ws.Range("M28").Value = -1009.125
# row 33
$ws.Range("H33").Value = 349.125
$ws.Range("I33").Value = 367.7143
$ws.Range("J33").Value = 219
$ws.Range("K33").Value = 367.7143
$ws.Range("L33").Value = 219
$ws.Range("M33").Value = -138.7143
$ws.Range("N33").Value = -677
# row 70
$ws.Range("H70").Value = 5400.385
$ws.Range("I70").Value = 4966.1113
$ws.Range("J70").Value = 5630.294
$ws.Range("K70").Value = 14898.3339
$ws.Range("L70").Value = 16890.882
$ws.Range("M70").Value = -14628.3339
$ws.Range("N70").Value = -17430.882
# row 73
$ws.Range("H73").Value = 5400.385
$ws.Range("I73").Value = 4966.1113
$ws.Range("J73").Value = 5630.294
$ws.Range("K73").Value = 14898.3339
$ws.Range("L73").Value = 16890.882
$ws.Range("M73").Value = -13962.3339
$ws.Range("N73").Value = -18762.882
# row 80
$ws.Range("H80").Value = 858.6667
$ws.Range("I80").Value = 423.6
$ws.Range("K80").Value = 1270.8
$ws.Range("M80").Value = -272.8000000000002
# row 83
$ws.Range("H83").Value = 858.6667
$ws.Range("I83").Value = 423.6
$ws.Range("K83").Value = 3812.4
$ws.Range("M83").Value = 1179.6
# row 86
$ws.Range("H86").Value = 17251.5
$ws.Range("I86").Value = 19499
$ws.Range("J86").Value = 15004
$ws.Range("K86").Value = 19499
$ws.Range("L86").Value = 15004
$ws.Range("M86").Value = -18376
$ws.Range("N86").Value = -17250
# row 88
$ws.Range("H88").Value = 1703.9231
$ws.Range("I88").Value = 2039.5
$ws.Range("K88").Value = 2039.5
$ws.Range("M88").Value = -1633.5
# row 89
$ws.Range("H89").Value = 17251.5
$ws.Range("I89").Value = 19499
$ws.Range("J89").Value = 15004
$ws.Range("K89").Value = 97495
$ws.Range("L89").Value = 75020
$ws.Range("M89").Value = -91879
$ws.Range("N89").Value = -86252
# row 91
$ws.Range("H91").Value = 1703.9231
$ws.Range("I91").Value = 2039.5
$ws.Range("K91").Value = 2039.5
$ws.Range("M91").Value = -635.5
# row 112
$ws.Range("H112").Value = 4153.6665
$ws.Range("J112").Value = 4189.643
$ws.Range("L112").Value = 12568.929
$ws.Range("N112").Value = -14784.929
# row 138
$ws.Range("H138").Value = 962.5454999999999
$ws.Range("I138").Value = 962.5454999999999
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2887.6365
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2252.3635
$ws.Range("N138").ClearContents()
# row 141
$ws.Range("H141").Value = 4531.2144
$ws.Range("I141").Value = 5030.5835
$ws.Range("K141").Value = 15091.7505
$ws.Range("M141").Value = -9911.750499999998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# row 55
$ws.Range("H55").Value = 13750
$ws.Range("I55").Value = 13750
$ws.Range("K55").Value = 13750
$ws.Range("M55").Value = -13435
# row 61
$ws.Range("H61").Value = 2430.9697
$ws.Range("I61").Value = 2299
$ws.Range("K61").Value = 2299
$ws.Range("M61").Value = -2087
# row 122
$ws.Range("H122").Value = 4988.5454
$ws.Range("I122").Value = 4041.6667
$ws.Range("J122").Value = 9249.5
$ws.Range("K122").Value = 12125.0001
$ws.Range("L122").Value = 27748.5
$ws.Range("M122").Value = -9675.000100000001
$ws.Range("N122").Value = -32648.5
# row 133
$ws.Range("H133").Value = 73261.664
$ws.Range("J133").Value = 73261.664
$ws.Range("L133").Value = 73261.664
$ws.Range("N133").Value = -78321.664
# row 136
$ws.Range("H136").Value = 2430.9697
$ws.Range("I136").Value = 2299
$ws.Range("K136").Value = 6897
$ws.Range("M136").Value = -4347

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 4434.7
$ws.Range("I86").Value = 2604.2666
$ws.Range("J86").Value = 9926
$ws.Range("K86").Value = 2604.2666
$ws.Range("L86").Value = 9926
$ws.Range("M86").Value = -1481.2666
$ws.Range("N86").Value = -12172
# row 89
$ws.Range("H89").Value = 4434.7
$ws.Range("I89").Value = 2604.2666
$ws.Range("J89").Value = 9926
$ws.Range("K89").Value = 13021.333
$ws.Range("L89").Value = 49630
$ws.Range("M89").Value = -7405.332999999999
$ws.Range("N89").Value = -60862
# row 99
$ws.Range("H99").Value = 2132.625
$ws.Range("I99").Value = 1609.4
$ws.Range("K99").Value = 1609.4
$ws.Range("M99").Value = -111.4000000000001
# row 107
$ws.Range("H107").Value = 1361.8948
$ws.Range("I107").Value = 1055.0625
$ws.Range("K107").Value = 1055.0625
$ws.Range("M107").Value = 864.9375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2777.9443
$ws.Range("I31").Value = 1362.2858
$ws.Range("K31").Value = 1362.2858
$ws.Range("M31").Value = -1067.2858
# row 34
$ws.Range("H34").Value = 2777.9443
$ws.Range("I34").Value = 1362.2858
$ws.Range("K34").Value = 1362.2858
$ws.Range("M34").Value = -1160.2858
# row 107
$ws.Range("H107").Value = 333.08334
$ws.Range("I107").Value = 272.45456
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 272.45456
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1647.54544
$ws.Range("N107").Value = -4840
# row 122
$ws.Range("H122").Value = 3980.111
$ws.Range("I122").Value = 4274.5
$ws.Range("J122").Value = 2949.75
$ws.Range("K122").Value = 12823.5
$ws.Range("L122").Value = 8849.25
$ws.Range("M122").Value = -10373.5
$ws.Range("N122").Value = -13749.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 121
$ws.Range("H121").Value = 667575.2
$ws.Range("J121").Value = 910035.5600000001
$ws.Range("L121").Value = 2730106.68
$ws.Range("N121").Value = -2732726.68
# row 134
$ws.Range("H134").Value = 1169.6666
$ws.Range("I134").Value = 1169.6666
$ws.Range("K134").Value = 3508.9998
$ws.Range("M134").Value = 1561.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 19
$ws.Range("H19").Value = 9000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 9000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 9000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -9576
# row 122
$ws.Range("H122").Value = 2581.6843
$ws.Range("I122").Value = 2517.5715
$ws.Range("J122").Value = 2761.2
$ws.Range("K122").Value = 7552.7145
$ws.Range("L122").Value = 8283.599999999999
$ws.Range("M122").Value = -5102.7145
$ws.Range("N122").Value = -13183.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 122
$ws.Range("H122").Value = 2890.1428
$ws.Range("I122").Value = 2843.2307
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8529.6921
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6079.6921
$ws.Range("N122").Value = -15400

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 19
$ws.Range("H19").Value = 4000
$ws.Range("J19").Value = 4000
$ws.Range("L19").Value = 4000
$ws.Range("N19").Value = -4348
# row 46
$ws.Range("H46").Value = 29998
$ws.Range("J46").Value = 29998
$ws.Range("L46").Value = 29998
$ws.Range("N46").Value = -30460
# row 93
$ws.Range("H93").Value = 57889
$ws.Range("J93").Value = 57889
$ws.Range("L93").Value = 57889
$ws.Range("N93").Value = -62881
# row 134
$ws.Range("H134").Value = 29998
$ws.Range("J134").Value = 29998
$ws.Range("L134").Value = 89994
$ws.Range("N134").Value = -95064
